$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.760.19'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.641.92'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.59%  '
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("D12").Value = '1.869.18'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '1.638.50'
$ws.Range("E13").Value = '  -0.80%  '
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.71'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").Value = '26.759.09'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("E22").Value = '  +6.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("E31").Value = '  +1.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.98%  '
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("D34").Value = '1.285.60'
$ws.Range("E34").Value = '  +0.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0177'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.818'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.807'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.64%  '
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("E43").Value = '  -2.61%  '
$ws.Range("D44").Value = '1.779.00'
$ws.Range("E44").Value = '  -0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.27'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0516'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0967'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.406'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.19%  '
